# Season-record columns: add Wins / Losses / Ties to the player table.
# (The previous scrape only pulled team statistics, not the season record.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (AC1, which carries
# the bold/centered/bordered header style) onto the three new header cells
# so they match the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row gets the team's season record repeated alongside them.
$lastRow = 39
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 83
    $ws.Cells.Item($r, 31).Value = 79
    $ws.Cells.Item($r, 32).Value = 0
}
